$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "215.39", "1.00", "19.75") are preserved as literal text, not converted
# to numbers by Excel's automatic type inference.
$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "25.781.06"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "1.635.77"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "215.39"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("E6").Value = "  -0.82%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("E9").Value = "  -1.08%  "

$ws.Range("D10").Value = "19.75"
$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +1.44%  "

$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").Value = "1.863.69"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").Value = "1.636.47"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  -0.54%  "

$ws.Range("D16").Value = "0.0₃0768"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").Value = "63.08"
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").Value = "25.810.57"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").Value = "4.45"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").Value = "193.24"
$ws.Range("E21").Value = "  -0.86%  "

$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("D23").Value = "6.34"
$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("E25").Value = "  +2.96%  "

$ws.Range("D26").Value = "142.61"
$ws.Range("E26").Value = "  +3.11%  "

$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("D28").Value = "6.96"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("D31").Value = "0.0495"
$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").Value = "3.25"
$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "0.903"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("D37").Value = "1.132.24"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  -1.84%  "

$ws.Range("D39").Value = "0.544"
$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("D42").Value = "5.56"
$ws.Range("E42").Value = "  +1.45%  "

$ws.Range("D43").Value = "100.54"
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("D45").Value = "1.773.31"
$ws.Range("E45").Value = "  -0.19%  "

$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -2.02%  "

$ws.Range("D47").Value = "55.30"

$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.41"
$ws.Range("E50").Value = "  +2.15%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.50"
$ws.Range("E51").Value = "  -3.20%  "

# Restore the default (Normal) style so the forced text number format
# does not linger on cells that did not have an explicit style before.
$rng.Style = "Normal"
